$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01306
$ws.Range("H2").Value = 0.03918
$ws.Range("I2").Value = 0.01051556883913273
$ws.Range("J2").Value = 0.01051556883913273
$ws.Range("Q2").Value = 0.0006233886266666667
$ws.Range("R2").Value = 0.00561049764
$ws.Range("S2").Value = 0.01051556883913273
$ws.Range("T2").Value = 0.01051556883913273

# Row 3
$ws.Range("I3").Value = 0.2290509900416114
$ws.Range("J3").Value = 0.2290509900416114
$ws.Range("S3").Value = 0.2290509900416114
$ws.Range("T3").Value = 0.2290509900416114

# Row 4
$ws.Range("G4").Value = 0.9444339999999999
$ws.Range("I4").Value = 0.7604334411192558
$ws.Range("J4").Value = 0.7604334411192558
$ws.Range("Q4").Value = 0.04508035331066666
$ws.Range("R4").Value = 0.4057231797959999
$ws.Range("S4").Value = 0.7604334411192558
$ws.Range("T4").Value = 0.7604334411192558
